# Update the mockup header row from Thai labels to the new English labels
# and refresh the active selection, per "add new school mockup data".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "index"
$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "name-surname"
$ws.Range("D1").Value = "class-room"
$ws.Range("E1").Value = "car-number"
$ws.Range("F1").Value = "father-phone"
$ws.Range("G1").Value = "mother-phone"
$ws.Range("H1").Value = "address"

$ws.Range("B4").Select()
